# fix: Criando tela de login
# Remove the "O Sistema deve identificar..." demand row from the "Demandas" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demandas")
$ws.Activate()

# Delete the entire row that holds the "password lockout" demand (row 4),
# shifting the following rows (and the table) up by one.
$ws.Rows("4:4").Delete()

# Restore cursor/selection position as saved in the workbook.
$ws.Range("C17").Select()

$wb.Save()
